# Refresh the cryptos snapshot: price (col D) and 1h volume-change (col E) figures,
# plus a row-44/45 coin swap (EnergySwap <-> FirstDigitalUSD), per the source update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.280.59'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '3.494.22'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''587.20'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '''134.10'
$ws.Range('E6').Value = '  +1.75%  '
$ws.Range('D7').Value = '3.494.08'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('D11').Value = '''7.18'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').Value = '4.087.96'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  +1.56%  '
$ws.Range('D15').Value = '''0.0000181'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '3.492.40'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '64.309.86'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '''25.38'
$ws.Range('E18').Value = '  -8.62%  '
$ws.Range('D19').Value = '''9.86'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('E21').Value = '  -5.58%  '
$ws.Range('D22').Value = '''388.31'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').Value = '''0.567'
$ws.Range('E23').Value = '  -1.95%  '
$ws.Range('D24').Value = '3.633.50'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '''74.33'
$ws.Range('E25').Value = '  +2.10%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('E31').Value = '  -4.96%  '
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('D33').Value = '''8.25'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('D34').Value = '3.515.51'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +3.82%  '
$ws.Range('D37').Value = '''23.48'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('D38').Value = '''5.24'
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('D39').Value = '''6.87'
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('E40').Value = '  -2.11%  '
$ws.Range('D41').Value = '''162.12'
$ws.Range('E41').Value = '  -2.91%  '
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''25.47'
$ws.Range('E45').Value = '  -5.05%  '
$ws.Range('D46').Value = '''41.82'
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('D50').Value = '2.470.52'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('E51').Value = '  -2.04%  '
